$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bdnf"
$ws.Range("C2").Value = "Sort1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03885866666666667
$ws.Range("H2").Value = 0.116576
$ws.Range("I2").Value = 0.01924839521029073
$ws.Range("J2").Value = 0.01924839521029073
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.358929333333333
$ws.Range("N2").Value = 13.076788
$ws.Range("O2").Value = 0.315486934607403
$ws.Range("P2").Value = 0.315486934607403
$ws.Range("Q2").Value = 0.1693821819875556
$ws.Range("R2").Value = 1.524439637888
$ws.Range("S2").Value = 0.00607261720100644
$ws.Range("T2").Value = 0.006072617201006441

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bdnf"
$ws.Range("C3").Value = "Sort1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03885866666666667
$ws.Range("H3").Value = 0.116576
$ws.Range("I3").Value = 0.01924839521029073
$ws.Range("J3").Value = 0.01924839521029073
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.781382333333333
$ws.Range("N3").Value = 5.344147
$ws.Range("O3").Value = 0.1289313977653647
$ws.Range("P3").Value = 0.1289313977653647
$ws.Range("Q3").Value = 0.06922214229688889
$ws.Range("R3").Value = 0.6229992806719999
$ws.Range("S3").Value = 0.002481722499202936
$ws.Range("T3").Value = 0.002481722499202936

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bdnf"
$ws.Range("C4").Value = "Sort1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03885866666666667
$ws.Range("H4").Value = 0.116576
$ws.Range("I4").Value = 0.01924839521029073
$ws.Range("J4").Value = 0.01924839521029073
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.676201333333334
$ws.Range("N4").Value = 23.028604
$ws.Range("O4").Value = 0.5555816676272323
$ws.Range("P4").Value = 0.5555816676272323
$ws.Range("Q4").Value = 0.2982869488782223
$ws.Range("R4").Value = 2.684582539904
$ws.Range("S4").Value = 0.01069405551008135
$ws.Range("T4").Value = 0.01069405551008135

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Bdnf"
$ws.Range("C5").Value = "Sort1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.979941666666667
$ws.Range("H5").Value = 5.939825
$ws.Range("I5").Value = 0.9807516047897092
$ws.Range("J5").Value = 0.9807516047897092
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.358929333333333
$ws.Range("N5").Value = 13.076788
$ws.Range("O5").Value = 0.315486934607403
$ws.Range("P5").Value = 0.315486934607403
$ws.Range("Q5").Value = 8.630425809122222
$ws.Range("R5").Value = 77.67383228210001
$ws.Range("S5").Value = 0.3094143174063965
$ws.Range("T5").Value = 0.3094143174063966

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Bdnf"
$ws.Range("C6").Value = "Sort1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.979941666666667
$ws.Range("H6").Value = 5.939825
$ws.Range("I6").Value = 0.9807516047897092
$ws.Range("J6").Value = 0.9807516047897092
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.781382333333333
$ws.Range("N6").Value = 5.344147
$ws.Range("O6").Value = 0.1289313977653647
$ws.Range("P6").Value = 0.1289313977653647
$ws.Range("Q6").Value = 3.527033106030555
$ws.Range("R6").Value = 31.743297954275
$ws.Range("S6").Value = 0.1264496752661618
$ws.Range("T6").Value = 0.1264496752661618

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Bdnf"
$ws.Range("C7").Value = "Sort1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.979941666666667
$ws.Range("H7").Value = 5.939825
$ws.Range("I7").Value = 0.9807516047897092
$ws.Range("J7").Value = 0.9807516047897092
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.676201333333334
$ws.Range("N7").Value = 23.028604
$ws.Range("O7").Value = 0.5555816676272323
$ws.Range("P7").Value = 0.5555816676272323
$ws.Range("Q7").Value = 15.19843086158889
$ws.Range("R7").Value = 136.7858777543
$ws.Range("S7").Value = 0.5448876121171509
$ws.Range("T7").Value = 0.5448876121171509
